$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 7081.5
$ws.Range("I74").Value = 3851.5
$ws.Range("J74").Value = 7889
$ws.Range("K74").Value = 3851.5
$ws.Range("L74").Value = 7889
$ws.Range("M74").Value = -2915.5
$ws.Range("N74").Value = -9761

$ws.Range("H77").Value = 7081.5
$ws.Range("I77").Value = 3851.5
$ws.Range("J77").Value = 7889
$ws.Range("K77").Value = 19257.5
$ws.Range("L77").Value = 39445
$ws.Range("M77").Value = -14577.5
$ws.Range("N77").Value = -48805

$ws.Range("H86").Value = 11225
$ws.Range("I86").Value = 17428.334
$ws.Range("J86").Value = 1920
$ws.Range("K86").Value = 17428.334
$ws.Range("L86").Value = 1920
$ws.Range("M86").Value = -16305.334
$ws.Range("N86").Value = -4166

$ws.Range("H89").Value = 11225
$ws.Range("I89").Value = 17428.334
$ws.Range("J89").Value = 1920
$ws.Range("K89").Value = 87141.67
$ws.Range("L89").Value = 9600
$ws.Range("M89").Value = -81525.67
$ws.Range("N89").Value = -20832

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10084.082
$ws.Range("I32").Value = 8662.201999999999
$ws.Range("J32").Value = 25902.5
$ws.Range("K32").Value = 8662.201999999999
$ws.Range("L32").Value = 25902.5
$ws.Range("M32").Value = -8375.201999999999
$ws.Range("N32").Value = -26476.5

$ws.Range("H61").Value = 1720.1455
$ws.Range("I61").Value = 1656.8966
$ws.Range("J61").Value = 1790.6923
$ws.Range("K61").Value = 1656.8966
$ws.Range("L61").Value = 1790.6923
$ws.Range("M61").Value = -1444.8966
$ws.Range("N61").Value = -2214.6923

$ws.Range("H122").Value = 11444.333
$ws.Range("I122").Value = 16666.5
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 49999.5
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -47549.5
$ws.Range("N122").Value = -7900

$ws.Range("H136").Value = 1720.1455
$ws.Range("I136").Value = 1656.8966
$ws.Range("J136").Value = 1790.6923
$ws.Range("K136").Value = 4970.6898
$ws.Range("L136").Value = 5372.0769
$ws.Range("M136").Value = -2420.6898
$ws.Range("N136").Value = -10472.0769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 4464.645
$ws.Range("I80").Value = 1002.1539
$ws.Range("J80").Value = 6965.3335
$ws.Range("K80").Value = 1002.1539
$ws.Range("L80").Value = 6965.3335
$ws.Range("M80").Value = -4.153900000000021
$ws.Range("N80").Value = -8961.333500000001

$ws.Range("H83").Value = 4464.645
$ws.Range("I83").Value = 1002.1539
$ws.Range("J83").Value = 6965.3335
$ws.Range("K83").Value = 5010.7695
$ws.Range("L83").Value = 34826.6675
$ws.Range("M83").Value = -18.76950000000033
$ws.Range("N83").Value = -44810.6675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3693.238
$ws.Range("I62").Value = 2881.5
$ws.Range("J62").Value = 4775.5557
$ws.Range("K62").Value = 2881.5
$ws.Range("L62").Value = 4775.5557
$ws.Range("M62").Value = -2257.5
$ws.Range("N62").Value = -6023.5557

$ws.Range("H65").Value = 3693.238
$ws.Range("I65").Value = 2881.5
$ws.Range("J65").Value = 4775.5557
$ws.Range("K65").Value = 14407.5
$ws.Range("L65").Value = 23877.7785
$ws.Range("M65").Value = -11287.5
$ws.Range("N65").Value = -30117.7785

$ws.Range("H122").Value = 30000868
$ws.Range("I122").Value = 37500780
$ws.Range("J122").Value = 1220
$ws.Range("K122").Value = 112502340
$ws.Range("L122").Value = 3660
$ws.Range("M122").Value = -112499890
$ws.Range("N122").Value = -8560

$ws.Range("H132").Value = 7940267.5
$ws.Range("I132").Value = 932.12
$ws.Range("J132").Value = 19615760
$ws.Range("K132").Value = 2796.36
$ws.Range("L132").Value = 58847280
$ws.Range("M132").Value = -266.3600000000001
$ws.Range("N132").Value = -58852340

$ws.Range("H134").Value = 1295.2307
$ws.Range("I134").Value = 1613.0883
$ws.Range("J134").Value = 694.8333
$ws.Range("K134").Value = 4839.2649
$ws.Range("L134").Value = 2084.4999
$ws.Range("M134").Value = -2304.2649
$ws.Range("N134").Value = -7154.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 28.368422
$ws.Range("I2").Value = 28.333334
$ws.Range("J2").Value = 28.4
$ws.Range("K2").Value = 170.000004
$ws.Range("L2").Value = 170.4
$ws.Range("M2").Value = -57.00000399999999
$ws.Range("N2").Value = -396.4

$ws.Range("H131").Value = 730.92
$ws.Range("J131").Value = 763.3555
$ws.Range("L131").Value = 2290.0665
$ws.Range("N131").Value = -12370.0665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 33341000
$ws.Range("I80").Value = 11499.5
$ws.Range("K80").Value = 11499.5
$ws.Range("M80").Value = -10501.5

$ws.Range("H83").Value = 33341000
$ws.Range("I83").Value = 11499.5
$ws.Range("K83").Value = 57497.5
$ws.Range("M83").Value = -52505.5

$ws.Range("H122").Value = 8200681.5
$ws.Range("I122").Value = 13162528
$ws.Range("J122").Value = 2849.087
$ws.Range("K122").Value = 39487584
$ws.Range("L122").Value = 8547.261
$ws.Range("M122").Value = -39485134
$ws.Range("N122").Value = -13447.261

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5209918
$ws.Range("I46").Value = 8333872.5
$ws.Range("J46").Value = 3326.6667
$ws.Range("K46").Value = 8333872.5
$ws.Range("L46").Value = 3326.6667
$ws.Range("M46").Value = -8333684.5
$ws.Range("N46").Value = -3702.6667

$ws.Range("H68").Value = 1583.3334
$ws.Range("I68").Value = 1642.8572
$ws.Range("K68").Value = 1642.8572
$ws.Range("M68").Value = -893.8571999999999

$ws.Range("H71").Value = 1583.3334
$ws.Range("I71").Value = 1642.8572
$ws.Range("K71").Value = 8214.286
$ws.Range("M71").Value = -4470.286

$ws.Range("H82").Value = 1325
$ws.Range("I82").Value = 1300
$ws.Range("J82").Value = 1400
$ws.Range("K82").Value = 1300
$ws.Range("L82").Value = 1400
$ws.Range("M82").Value = -939
$ws.Range("N82").Value = -2122

$ws.Range("H85").Value = 1325
$ws.Range("I85").Value = 1300
$ws.Range("J85").Value = 1400
$ws.Range("K85").Value = 1300
$ws.Range("L85").Value = 1400
$ws.Range("M85").Value = -52
$ws.Range("N85").Value = -3896

$ws.Range("H93").Value = 1189.0416
$ws.Range("I93").Value = 1152
$ws.Range("J93").Value = 1279
$ws.Range("K93").Value = 1152
$ws.Range("L93").Value = 1279
$ws.Range("M93").Value = 96
$ws.Range("N93").Value = -3775

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1739.5428
$ws.Range("I122").Value = 1079.9615
$ws.Range("K122").Value = 3239.8845
$ws.Range("M122").Value = -789.8844999999997
